$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F7").Value = 246
$ws.Range("F9").Value = 342
$ws.Range("F12").Value = 684
$ws.Range("F14").Value = 507
$ws.Range("F15").Value = 142
$ws.Range("F19").Value = 2620
$ws.Range("F27").Value = 591
$ws.Range("F28").Value = 984

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F14").Value = 610
$ws.Range("F17").Value = 987
$ws.Range("F26").Value = 3930
$ws.Range("F33").Value = 163
$ws.Range("F35").Value = 11

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 1781
$ws.Range("F5").Value = 2455
$ws.Range("F6").Value = 1042
$ws.Range("F9").Value = 1321

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 1781
$ws.Range("F4").Value = 2455
$ws.Range("F6").Value = 1042
$ws.Range("F7").Value = 1321
$ws.Range("F13").Value = 246
$ws.Range("F15").Value = 342
$ws.Range("F17").Value = 684
$ws.Range("F20").Value = 507
$ws.Range("F24").Value = 2620
$ws.Range("F31").Value = 591
$ws.Range("F32").Value = 984
$ws.Range("F33").Value = 610
$ws.Range("F34").Value = 610
$ws.Range("F49").Value = 163
